$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Winter collection" colors appended to the COR / ID table (rows 19-23).
# Column A keeps using text (shared strings), column B now holds real numbers
# (18-22) instead of the old " N" text labels used in rows 2-18.
$ws.Range("A19").Value = "SALMÃO"
$ws.Range("B19").Value = 18
$ws.Range("A20").Value = "ROSÊ"
$ws.Range("B20").Value = 19
$ws.Range("A21").Value = "FOSSIL"
$ws.Range("B21").Value = 20
$ws.Range("A22").Value = "TIJOLO"
$ws.Range("B22").Value = 21
$ws.Range("A23").Value = "CREME"
$ws.Range("B23").Value = 22

# Apply a single combined style (integer number format + left alignment) to
# the whole column B range. Building the format on a scratch cell first and
# then pasting just the formats avoids creating one throw-away cell style per
# property set on the destination range.
$temp = $ws.Range("Z1")
$temp.NumberFormat = "0"
$temp.HorizontalAlignment = -4131
$colB = $ws.Range("B1:B23")
$temp.Copy()
$colB.PasteSpecial(-4122)
$temp.Clear()

# Match the author's final selection recorded in the saved view state.
$ws.Range("M5").Select() | Out-Null
